$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Map of row label -> array of new values for columns 2..5 (Distance chi2, Distance p, Urb chi2, Urb p)
$updates = @{
    "Danaus plexippus abundance"       = @("1.545", "0.214", "1.987", "0.159")
    "Liriomyza asclepiadis abundance"  = @("0.092", "0.762", "0.666", "0.414")
    "Labidomera clivicollis abundance" = @("0.017", "0.898", "0.198", "0.656")
    "SLA"                              = @("0.100", "0.752", "0.135", "0.714")
    "Height before flowering"          = @("1.275", "0.259", "0.304", "0.581")
}

for ($ri = 1; $ri -le $t.Rows.Count; $ri++) {
    $row = $t.Rows.Item($ri)
    $label = $row.Cells.Item(1).Range.Text

    foreach ($key in $updates.Keys) {
        if ($label -like "$key*") {
            $newVals = $updates[$key]
            for ($ci = 2; $ci -le 5; $ci++) {
                $row.Cells.Item($ci).Range.Text = $newVals[$ci - 2]
            }
        }
    }
}
